$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 87

$ws.Range("E10").Value = 400
$ws.Range("F10").Value = 190
$ws.Range("H10").Value = 190

$ws.Range("E12").Value = 387

$ws.Range("E23").Value = 169

$ws.Range("E34").Value = 173

$ws.Range("E42").Value = 290

$ws.Range("E48").Value = 171
